$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 185, shifting existing rows 185:290 down to 186:291
$ws.Rows.Item(185).Insert()

# Populate the newly inserted row 185 with the new record
$ws.Cells.Item(185, 1).Value  = 10
$ws.Cells.Item(185, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(185, 3).Value  = "La Araucanía"
$ws.Cells.Item(185, 4).Value  = 44455
$ws.Cells.Item(185, 5).Value  = 9
$ws.Cells.Item(185, 6).Value  = 100112043
$ws.Cells.Item(185, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(185, 8).Value  = "Sin especificar"
$ws.Cells.Item(185, 9).Value  = "Primera"
$ws.Cells.Item(185, 10).Value = 60
$ws.Cells.Item(185, 11).Value = 16000
$ws.Cells.Item(185, 12).Value = 17000
$ws.Cells.Item(185, 13).Value = 16500
$ws.Cells.Item(185, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(185, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(185, 16).Value = 275
$ws.Cells.Item(185, 17).Value = 60
$ws.Cells.Item(185, 18).Value = "Hortaliza"
